$d = $word.ActiveDocument

$pairs = @(
    @("12×89=1068", "16×21=336"),
    @("66×55=3630", "41×90=3690"),
    @("50×54=2700", "11×93=1023"),
    @("58×99=5742", "75×14=1050"),
    @("11×88=968",  "16×19=304"),
    @("36×84=3024", "61×19=1159"),
    @("52×11=572",  "59×51=3009"),
    @("93×19=1767", "11×21=231"),
    @("86×23=1978", "15×12=180"),
    @("94×13=1222", "59×27=1593"),
    @("58×71=4118", "42×63=2646"),
    @("72×66=4752", "95×25=2375"),
    @("59×74=4366", "50×50=2500"),
    @("80×26=2080", "95×65=6175"),
    @("73×84=6132", "66×88=5808"),
    @("80×57=4560", "69×23=1587"),
    @("49×54=2646", "81×66=5346"),
    @("43×96=4128", "72×87=6264"),
    @("90×53=4770", "28×86=2408"),
    @("60×54=3240", "82×68=5576"),
    @("17×73=1241", "25×30=750"),
    @("61×56=3416", "74×58=4292"),
    @("81×38=3078", "36×13=468"),
    @("58×60=3480", "92×34=3128"),
    @("35×84=2940", "12×95=1140")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
